$d = $word.ActiveDocument

$find = " is allowable for type parameter constrained to "
$replace = " is allowable for type parameter categorized as "

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
